$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 4.920600000000002
$ws.Range("A10").Value = -21.88529999999998
$ws.Range("A12").Value = -21.56630000000001
$ws.Range("B15").Value = 4.695799999999996
$ws.Range("A18").Value = -22.38650000000002
$ws.Range("B20").Value = 9.487599999999993
$ws.Range("B29").Value = 4.932100000000003
$ws.Range("B30").Value = 5.761600000000002
$ws.Range("B31").Value = 5.127900000000002
$ws.Range("A37").Value = -19.2103
$ws.Range("B40").Value = 9.236899999999993
$ws.Range("A55").Value = -21.727
$ws.Range("A68").Value = -21.49690000000001
$ws.Range("B68").Value = 4.539000000000001
$ws.Range("B76").Value = 5.693899999999997
$ws.Range("A77").Value = -20.464
$ws.Range("A78").Value = -19.91559999999998
$ws.Range("B87").Value = 5.074399999999996
$ws.Range("B88").Value = 4.510699999999998
$ws.Range("B96").Value = 5.423800000000005
$ws.Range("B98").Value = 5.628299999999999
$ws.Range("B101").Value = 9.505999999999995
$ws.Range("B102").Value = 8.606600000000006
